$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated calibration data for rows 1-90 (columns A-D)
$arr = New-Object 'object[,]' 90,4
$arr[0,0] = 926.125
$arr[0,1] = 1014.961
$arr[0,2] = 2500.8
$arr[0,3] = 0
$arr[1,0] = 963.603
$arr[1,1] = 1788.024
$arr[1,2] = 2500.8
$arr[1,3] = -15.9
$arr[2,0] = 956.566
$arr[2,1] = 1692.296
$arr[2,2] = 2500.8
$arr[2,3] = -14.08
$arr[3,0] = 949.801
$arr[3,1] = 1588.279
$arr[3,2] = 2500.8
$arr[3,3] = -12.08
$arr[4,0] = 943.714
$arr[4,1] = 1486.918
$arr[4,2] = 2500.8
$arr[4,3] = -10.08
$arr[5,0] = 938.557
$arr[5,1] = 1385.821
$arr[5,2] = 2500.8
$arr[5,3] = -8.08
$arr[6,0] = 934.062
$arr[6,1] = 1284.03
$arr[6,2] = 2500.8
$arr[6,3] = -6.08
$arr[7,0] = 930.527
$arr[7,1] = 1183.591
$arr[7,2] = 2500.8
$arr[7,3] = -4.08
$arr[8,0] = 927.904
$arr[8,1] = 1082.109
$arr[8,2] = 2500.8
$arr[8,3] = -2.08
$arr[9,0] = 925.915
$arr[9,1] = 982.938
$arr[9,2] = 2500.8
$arr[9,3] = -0.08
$arr[10,0] = 924.533
$arr[10,1] = 881.414
$arr[10,2] = 2500.8
$arr[10,3] = 1.92
$arr[11,0] = 924.0
$arr[11,1] = 780.127
$arr[11,2] = 2500.8
$arr[11,3] = 3.92
$arr[12,0] = 924.457
$arr[12,1] = 679.393
$arr[12,2] = 2500.8
$arr[12,3] = 5.92
$arr[13,0] = 925.594
$arr[13,1] = 579.15
$arr[13,2] = 2500.8
$arr[13,3] = 7.92
$arr[14,0] = 927.547
$arr[14,1] = 477.93
$arr[14,2] = 2500.8
$arr[14,3] = 9.95
$arr[15,0] = 930.241
$arr[15,1] = 376.182
$arr[15,2] = 2500.8
$arr[15,3] = 11.92
$arr[16,0] = 933.617
$arr[16,1] = 272.609
$arr[16,2] = 2500.8
$arr[16,3] = 13.97
$arr[17,0] = 937.846
$arr[17,1] = 168.371
$arr[17,2] = 2500.8
$arr[17,3] = 15.97
$arr[18,0] = 907.354
$arr[18,1] = 1015.435
$arr[18,2] = 2600.5
$arr[18,3] = 0
$arr[19,0] = 946.873
$arr[19,1] = 1823.058
$arr[19,2] = 2600.5
$arr[19,3] = -15.9
$arr[20,0] = 939.825
$arr[20,1] = 1728.217
$arr[20,2] = 2600.5
$arr[20,3] = -14.07
$arr[21,0] = 932.869
$arr[21,1] = 1625.484
$arr[21,2] = 2600.5
$arr[21,3] = -12.1
$arr[22,0] = 926.702
$arr[22,1] = 1522.644
$arr[22,2] = 2600.5
$arr[22,3] = -10.1
$arr[23,0] = 921.408
$arr[23,1] = 1421.28
$arr[23,2] = 2600.5
$arr[23,3] = -8.07
$arr[24,0] = 916.714
$arr[24,1] = 1321.503
$arr[24,2] = 2600.5
$arr[24,3] = -6.1
$arr[25,0] = 912.979
$arr[25,1] = 1220.499
$arr[25,2] = 2600.5
$arr[25,3] = -4.1
$arr[26,0] = 909.969
$arr[26,1] = 1119.784
$arr[26,2] = 2600.5
$arr[26,3] = -2.1
$arr[27,0] = 907.598
$arr[27,1] = 1019.811
$arr[27,2] = 2600.5
$arr[27,3] = -0.1
$arr[28,0] = 905.96
$arr[28,1] = 919.26
$arr[28,2] = 2600.5
$arr[28,3] = 1.9
$arr[29,0] = 905.092
$arr[29,1] = 817.714
$arr[29,2] = 2600.5
$arr[29,3] = 3.9
$arr[30,0] = 905.062
$arr[30,1] = 717.581
$arr[30,2] = 2600.5
$arr[30,3] = 5.9
$arr[31,0] = 906.0
$arr[31,1] = 616.054
$arr[31,2] = 2600.5
$arr[31,3] = 7.93
$arr[32,0] = 907.575
$arr[32,1] = 515.12
$arr[32,2] = 2600.5
$arr[32,3] = 9.93
$arr[33,0] = 909.934
$arr[33,1] = 414.106
$arr[33,2] = 2600.5
$arr[33,3] = 11.95
$arr[34,0] = 912.587
$arr[34,1] = 312.005
$arr[34,2] = 2600.5
$arr[34,3] = 13.95
$arr[35,0] = 916.402
$arr[35,1] = 207.761
$arr[35,2] = 2600.5
$arr[35,3] = 15.95
$arr[36,0] = 890.011
$arr[36,1] = 1015.314
$arr[36,2] = 2700.2
$arr[36,3] = 0
$arr[37,0] = 928.773
$arr[37,1] = 1821.647
$arr[37,2] = 2700.2
$arr[37,3] = -15.9
$arr[38,0] = 921.797
$arr[38,1] = 1727.053
$arr[38,2] = 2700.2
$arr[38,3] = -14.1
$arr[39,0] = 914.905
$arr[39,1] = 1622.457
$arr[39,2] = 2700.2
$arr[39,3] = -12.07
$arr[40,0] = 908.956
$arr[40,1] = 1522.005
$arr[40,2] = 2700.2
$arr[40,3] = -10.1
$arr[41,0] = 903.644
$arr[41,1] = 1420.732
$arr[41,2] = 2700.2
$arr[41,3] = -8.1
$arr[42,0] = 899.115
$arr[42,1] = 1320.275
$arr[42,2] = 2700.2
$arr[42,3] = -6.1
$arr[43,0] = 895.269
$arr[43,1] = 1219.012
$arr[43,2] = 2700.2
$arr[43,3] = -4.1
$arr[44,0] = 892.248
$arr[44,1] = 1118.978
$arr[44,2] = 2700.2
$arr[44,3] = -2.1
$arr[45,0] = 890.057
$arr[45,1] = 1018.561
$arr[45,2] = 2700.2
$arr[45,3] = -0.1
$arr[46,0] = 888.5
$arr[46,1] = 917.869
$arr[46,2] = 2700.2
$arr[46,3] = 1.9
$arr[47,0] = 887.896
$arr[47,1] = 817.547
$arr[47,2] = 2700.2
$arr[47,3] = 3.9
$arr[48,0] = 887.9
$arr[48,1] = 716.956
$arr[48,2] = 2700.2
$arr[48,3] = 5.9
$arr[49,0] = 888.483
$arr[49,1] = 615.564
$arr[49,2] = 2700.2
$arr[49,3] = 7.93
$arr[50,0] = 889.987
$arr[50,1] = 514.627
$arr[50,2] = 2700.2
$arr[50,3] = 9.93
$arr[51,0] = 891.985
$arr[51,1] = 413.658
$arr[51,2] = 2700.2
$arr[51,3] = 11.95
$arr[52,0] = 894.676
$arr[52,1] = 311.698
$arr[52,2] = 2700.2
$arr[52,3] = 13.95
$arr[53,0] = 898.216
$arr[53,1] = 207.544
$arr[53,2] = 2700.2
$arr[53,3] = 15.97
$arr[54,0] = 873.936
$arr[54,1] = 1015.219
$arr[54,2] = 2799.9
$arr[54,3] = 0
$arr[55,0] = 912.001
$arr[55,1] = 1821.643
$arr[55,2] = 2799.9
$arr[55,3] = -15.9
$arr[56,0] = 905.137
$arr[56,1] = 1726.567
$arr[56,2] = 2799.9
$arr[56,3] = -14.1
$arr[57,0] = 898.38
$arr[57,1] = 1623.004
$arr[57,2] = 2799.9
$arr[57,3] = -12.07
$arr[58,0] = 892.615
$arr[58,1] = 1522.052
$arr[58,2] = 2799.9
$arr[58,3] = -10.1
$arr[59,0] = 887.461
$arr[59,1] = 1420.51
$arr[59,2] = 2799.9
$arr[59,3] = -8.1
$arr[60,0] = 883.032
$arr[60,1] = 1320.364
$arr[60,2] = 2799.9
$arr[60,3] = -6.1
$arr[61,0] = 879.094
$arr[61,1] = 1219.396
$arr[61,2] = 2799.9
$arr[61,3] = -4.1
$arr[62,0] = 876.085
$arr[62,1] = 1118.732
$arr[62,2] = 2799.9
$arr[62,3] = -2.1
$arr[63,0] = 873.997
$arr[63,1] = 1019.467
$arr[63,2] = 2799.9
$arr[63,3] = -0.1
$arr[64,0] = 872.439
$arr[64,1] = 918.688
$arr[64,2] = 2799.9
$arr[64,3] = 1.9
$arr[65,0] = 871.808
$arr[65,1] = 818.596
$arr[65,2] = 2799.9
$arr[65,3] = 3.9
$arr[66,0] = 871.693
$arr[66,1] = 717.976
$arr[66,2] = 2799.9
$arr[66,3] = 5.9
$arr[67,0] = 872.037
$arr[67,1] = 616.9
$arr[67,2] = 2799.9
$arr[67,3] = 7.93
$arr[68,0] = 873.604
$arr[68,1] = 515.733
$arr[68,2] = 2799.9
$arr[68,3] = 9.93
$arr[69,0] = 875.534
$arr[69,1] = 414.451
$arr[69,2] = 2799.9
$arr[69,3] = 11.95
$arr[70,0] = 878.403
$arr[70,1] = 312.53
$arr[70,2] = 2799.9
$arr[70,3] = 13.95
$arr[71,0] = 881.417
$arr[71,1] = 209.196
$arr[71,2] = 2799.9
$arr[71,3] = 15.97
$arr[72,0] = 858.913
$arr[72,1] = 1015.312
$arr[72,2] = 2899.6
$arr[72,3] = 0
$arr[73,0] = 896.037
$arr[73,1] = 1819.271
$arr[73,2] = 2899.6
$arr[73,3] = -15.9
$arr[74,0] = 889.629
$arr[74,1] = 1724.83
$arr[74,2] = 2899.6
$arr[74,3] = -14.1
$arr[75,0] = 883.067
$arr[75,1] = 1621.659
$arr[75,2] = 2899.6
$arr[75,3] = -12.08
$arr[76,0] = 877.47
$arr[76,1] = 1520.89
$arr[76,2] = 2899.6
$arr[76,3] = -10.1
$arr[77,0] = 872.24
$arr[77,1] = 1419.481
$arr[77,2] = 2899.6
$arr[77,3] = -8.1
$arr[78,0] = 867.996
$arr[78,1] = 1319.373
$arr[78,2] = 2899.6
$arr[78,3] = -6.1
$arr[79,0] = 864.066
$arr[79,1] = 1218.388
$arr[79,2] = 2899.6
$arr[79,3] = -4.1
$arr[80,0] = 861.166
$arr[80,1] = 1117.871
$arr[80,2] = 2899.6
$arr[80,3] = -2.1
$arr[81,0] = 858.991
$arr[81,1] = 1018.205
$arr[81,2] = 2899.6
$arr[81,3] = -0.1
$arr[82,0] = 857.5
$arr[82,1] = 918.019
$arr[82,2] = 2899.6
$arr[82,3] = 1.9
$arr[83,0] = 856.77
$arr[83,1] = 817.527
$arr[83,2] = 2899.6
$arr[83,3] = 3.9
$arr[84,0] = 856.749
$arr[84,1] = 716.94
$arr[84,2] = 2899.6
$arr[84,3] = 5.9
$arr[85,0] = 857.004
$arr[85,1] = 616.168
$arr[85,2] = 2899.6
$arr[85,3] = 7.92
$arr[86,0] = 858.427
$arr[86,1] = 515.626
$arr[86,2] = 2899.6
$arr[86,3] = 9.92
$arr[87,0] = 860.193
$arr[87,1] = 413.87
$arr[87,2] = 2899.6
$arr[87,3] = 11.95
$arr[88,0] = 862.62
$arr[88,1] = 312.123
$arr[88,2] = 2899.6
$arr[88,3] = 13.95
$arr[89,0] = 865.836
$arr[89,1] = 208.521
$arr[89,2] = 2899.6
$arr[89,3] = 15.97
$ws.Range("A1:D90").Value = $arr

# Remove the old trailing data block (rows 91-108), keeping formatting
$ws.Range("A91:D108").ClearContents()

# Select the active data range, matching the saved view state
$ws.Range("A1:D90").Select()
